$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

$ws.Range("B2").Value = 0.79263823653846
$ws.Range("C2").Value = 0.705740642307693
$ws.Range("D2").Value = 0.7128991173076918
$ws.Range("E2").Value = 1.51125994326923
$ws.Range("F2").Value = 0.7189450956730764

$ws.Range("B3").Value = 1.928472638942304
$ws.Range("C3").Value = 2.100594170192307
$ws.Range("D3").Value = 1.717727036057692
$ws.Range("E3").Value = 2.385732178125001
$ws.Range("F3").Value = 1.727314515865384

$ws.Range("B4").Value = 0.7579752721153847
$ws.Range("C4").Value = 0.9449168495192313
$ws.Range("D4").Value = 0.9685790293269226
$ws.Range("E4").Value = 1.310371479807694
$ws.Range("F4").Value = 0.9139563562500003

$ws.Range("B5").Value = 2.464334802403844
$ws.Range("C5").Value = 1.899966336538461
$ws.Range("D5").Value = 1.607334240384615
$ws.Range("E5").Value = 2.646548055769231
$ws.Range("F5").Value = 1.577224537500001
